# Apply "semana 42 de 2025" updates to the poisson data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D3"  = 5
    "E3"  = 0
    "D4"  = 0
    "E4"  = 1
    "C5"  = 7
    "D5"  = 4
    "E5"  = 0.09
    "C6"  = 2
    "D6"  = 3
    "E6"  = 0.18
    "C7"  = 1
    "D7"  = 1
    "E7"  = 0.37
    "C9"  = 46
    "D9"  = 51
    "C11" = 1
    "D11" = 1
    "E11" = 0.37
    "C12" = 4
    "D12" = 2
    "E12" = 0.15
    "C16" = 1
    "E16" = 0.37
    "C17" = 13
    "D17" = 17
    "E17" = 0.05
    "C18" = 1
    "E18" = 0.37
    "C19" = 11
    "D19" = 6
    "E19" = 0.04
    "C20" = 1
    "E20" = 0.37
    "D22" = 1
    "E22" = 0.27
    "C26" = 1
    "D26" = 1
    "E26" = 0.37
    "D30" = 2
    "E30" = 0
    "C31" = 2
    "D31" = 2
    "E31" = 0.27
    "C33" = 5
    "D33" = 2
    "E33" = 0.08
    "C34" = 8
    "E34" = 0.01
    "C35" = 9
    "D35" = 11
    "E35" = 0.1
    "D36" = 0
    "E36" = 1
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
